$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('I2').Value = 7281
$ws.Range('J2').Value = 6259
$ws.Range('J3').Value = 6664
$ws.Range('J4').Value = 1447
$ws.Range('J5').Value = 511
$ws.Range('J6').Value = 8731
$ws.Range('I7').Value = 26234
$ws.Range('J7').Value = 23612

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('J2').Value = 58
$ws.Range('J6').Value = 225
$ws.Range('J7').Value = 338

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('J2').Value = 405
$ws.Range('J3').Value = 450
$ws.Range('J5').Value = 37
$ws.Range('J6').Value = 522
$ws.Range('J7').Value = 1492

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('J2').Value = 253
$ws.Range('J3').Value = 353
$ws.Range('J4').Value = 44
$ws.Range('J6').Value = 378
$ws.Range('J7').Value = 1072

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range('J2').Value = 115
$ws.Range('J7').Value = 337

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('J2').Value = 211
$ws.Range('J6').Value = 213
$ws.Range('J7').Value = 721

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('J2').Value = 187
$ws.Range('J6').Value = 178
$ws.Range('J7').Value = 689
$ws.Range('J8').Value = 1492
$ws.Range('J10').Value = 172
$ws.Range('J11').Value = 388
$ws.Range('J19').Value = 696
$ws.Range('J23').Value = 220
$ws.Range('J27').Value = 145
$ws.Range('G29').Value = 1801
$ws.Range('J29').Value = 1299
$ws.Range('J31').Value = 219
$ws.Range('J33').Value = 1072
$ws.Range('J36').Value = 318
$ws.Range('J37').Value = 721
$ws.Range('J41').Value = 159
$ws.Range('J42').Value = 1016
$ws.Range('J43').Value = 200
$ws.Range('J50').Value = 142
$ws.Range('J51').Value = 294
$ws.Range('J52').Value = 591
$ws.Range('J53').Value = 338
$ws.Range('J54').Value = 449
$ws.Range('G63').Value = 275
$ws.Range('I63').Value = 249
$ws.Range('J63').Value = 82
$ws.Range('J67').Value = 892
$ws.Range('J68').Value = 50
$ws.Range('J72').Value = 93
$ws.Range('J73').Value = 229
$ws.Range('J74').Value = 27
$ws.Range('J76').Value = 358
$ws.Range('J77').Value = 176
$ws.Range('J79').Value = 665
$ws.Range('J85').Value = 982
$ws.Range('J87').Value = 77
$ws.Range('J88').Value = 248
$ws.Range('J90').Value = 255
$ws.Range('J91').Value = 270
$ws.Range('J94').Value = 249
$ws.Range('J95').Value = 337
$ws.Range('J96').Value = 262
$ws.Range('J97').Value = 208
$ws.Range('J98').Value = 178
$ws.Range('I101').Value = 26234
$ws.Range('J101').Value = 23612

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('J2').Value = 80
$ws.Range('J7').Value = 219

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('J2').Value = 224
$ws.Range('J6').Value = 242
$ws.Range('J7').Value = 892

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('J3').Value = 90
$ws.Range('J7').Value = 449

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('J2').Value = 397
$ws.Range('G3').Value = 750
$ws.Range('J5').Value = 51
$ws.Range('J6').Value = 328
$ws.Range('G7').Value = 1801
$ws.Range('J7').Value = 1299

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('J6').Value = 270
$ws.Range('J7').Value = 696

$ws = $wb.Worksheets.Item('River North')
$ws.Range('J3').Value = 71
$ws.Range('J7').Value = 358

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range('J3').Value = 43
$ws.Range('J7').Value = 178

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range('J6').Value = 93
$ws.Range('J7').Value = 159

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('J2').Value = 216
$ws.Range('J6').Value = 540
$ws.Range('J7').Value = 1016

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range('J6').Value = 94
$ws.Range('J7').Value = 172

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range('J3').Value = 75
$ws.Range('J7').Value = 220

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('J6').Value = 92
$ws.Range('J7').Value = 262

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('J3').Value = 112
$ws.Range('J6').Value = 67
$ws.Range('J7').Value = 270

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('J6').Value = 198
$ws.Range('J7').Value = 665

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('J4').Value = 14
$ws.Range('J7').Value = 318

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('J2').Value = 214
$ws.Range('J7').Value = 689

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('J6').Value = 138
$ws.Range('J7').Value = 249

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range('J2').Value = 32
$ws.Range('J6').Value = 112
$ws.Range('J7').Value = 178

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range('J6').Value = 47
$ws.Range('J7').Value = 142

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('J2').Value = 112
$ws.Range('J4').Value = 24
$ws.Range('J6').Value = 172
$ws.Range('J7').Value = 388

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('J6').Value = 79
$ws.Range('J7').Value = 229

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range('J3').Value = 47
$ws.Range('J7').Value = 187

$ws = $wb.Worksheets.Item('West Town')
$ws.Range('J6').Value = 145
$ws.Range('J7').Value = 208

$ws = $wb.Worksheets.Item('United Center')
$ws.Range('J3').Value = 64
$ws.Range('J6').Value = 118
$ws.Range('J7').Value = 248

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range('J6').Value = 52
$ws.Range('J7').Value = 145

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range('J2').Value = 91
$ws.Range('J7').Value = 255

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('J3').Value = 77
$ws.Range('J7').Value = 294

$ws = $wb.Worksheets.Item('North Park')
$ws.Range('J4').Value = 4
$ws.Range('J6').Value = 12
$ws.Range('J7').Value = 50

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range('J4').Value = 19
$ws.Range('J7').Value = 200

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('J2').Value = 261
$ws.Range('J3').Value = 349
$ws.Range('J6').Value = 283
$ws.Range('J7').Value = 982

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range('J4').Value = 9
$ws.Range('J7').Value = 93

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range('J2').Value = 67
$ws.Range('J3').Value = 58
$ws.Range('J7').Value = 176

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('J4').Value = 21
$ws.Range('J7').Value = 591

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Range('J6').Value = 52
$ws.Range('J7').Value = 77

$ws = $wb.Worksheets.Item('Printers Row')
$ws.Range('J3').Value = 9
$ws.Range('J7').Value = 27
